# Update the NATMI ligand/receptor TPM table (Angpt1 -> Itgb1) with the new
# TPM-based numbers. A new sending/target cluster, "Resolving-Mac", is added
# to the existing 3 clusters (ECs, FAPs, MuSCs) which turns the previous
# 3x4 (sender x target) block into a full 4x4 block -> 16 data rows instead
# of 12, i.e. rows 2..17 instead of rows 2..13.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$clusters = @("ECs", "FAPs", "MuSCs", "Resolving-Mac")

# Columns E..J ("Ligand-expressing cells" .. "Ligand derived specificity of
# total expression value") depend only on the sending cluster.
$bySender = @{
    "ECs"           = @(2, 0.6666666666666666, 0.1786683333333333, 0.536005, 0.005700931134181372, 0.005700931134181372)
    "FAPs"          = @(3, 1, 22.14783133333333, 66.443494, 0.7066907652137446, 0.7066907652137446)
    "MuSCs"         = @(3, 1, 8.967326666666667, 26.90198, 0.2861285535640985, 0.2861285535640984)
    "Resolving-Mac" = @(2, 0.6666666666666666, 0.04637566666666667, 0.139127, 0.001479750087975396, 0.001479750087975395)
}

# Columns K..P ("Receptor-expressing cells" .. "Receptor derived specificity
# of total expression value") depend only on the target cluster.
$byTarget = @{
    "ECs"           = @(3, 1, 145.7007446666667, 437.1022340000001, 0.2865937750105843, 0.2865937750105843)
    "FAPs"          = @(3, 1, 168.7997026666667, 506.3991080000001, 0.3320294904365841, 0.3320294904365841)
    "MuSCs"         = @(3, 1, 128.1261546666667, 384.378464, 0.2520245069956105, 0.2520245069956105)
    "Resolving-Mac" = @(3, 1, 65.761079, 197.283237, 0.1293522275572212, 0.1293522275572212)
}

# Columns Q..T ("Edge average/total expression weight/derived specificity")
# are specific to each (sender, target) pair.
$byPair = @{
    "ECs|ECs"                     = @(26.03210921501889, 234.28898293517, 0.001633851374820411, 0.001633851374820411)
    "ECs|FAPs"                    = @(30.15916154261556, 271.43245388354, 0.001892877259496299, 0.001892877259496299)
    "ECs|MuSCs"                   = @(22.89208651070222, 206.02877859632, 0.001436774358507987, 0.001436774358507987)
    "ECs|Resolving-Mac"           = @(11.74942238313166, 105.744801448185, 0.0007374281413566761, 0.0007374281413566761)
    "FAPs|ECs"                    = @(3226.9555180184, 29042.5996621656, 0.2025331741677255, 0.2025331741677255)
    "FAPs|FAPs"                   = @(3738.547343778151, 33646.92609400336, 0.2346421746701594, 0.2346421746701594)
    "FAPs|MuSCs"                  = @(2837.716462945913, 25539.44816651322, 0.1781033917013447, 0.1781033917013447)
    "FAPs|Resolving-Mac"          = @(1456.465285990008, 13108.18757391008, 0.09141202467451509, 0.09141202467451509)
    "MuSCs|ECs"                   = @(1306.546173002591, 11758.91555702332, 0.08200266230425315, 0.08200266230425313)
    "MuSCs|FAPs"                  = @(1513.682075048205, 13623.13867543384, 0.09500311783924449, 0.09500311783924448)
    "MuSCs|MuSCs"                 = @(1148.949083439858, 10340.54175095872, 0.07211140764935904, 0.07211140764935903)
    "MuSCs|Resolving-Mac"         = @(589.7010773454733, 5307.30969610926, 0.03701136577124183, 0.03701136577124182)
    "Resolving-Mac|ECs"           = @(6.756969167746446, 60.81272250971801, 0.0004240871637851128, 0.0004240871637851127)
    "Resolving-Mac|FAPs"          = @(7.828198744301779, 70.45378869871601, 0.0004913206676839611, 0.000491320667683961)
    "Resolving-Mac|MuSCs"         = @(5.941935840103112, 53.477422560928, 0.0003729332863987103, 0.0003729332863987102)
    "Resolving-Mac|Resolving-Mac" = @(3.049713879344333, 27.447424914099, 0.0001914089701076115, 0.0001914089701076115)
}

$row = 2
foreach ($sender in $clusters) {
    foreach ($target in $clusters) {
        $ws.Cells.Item($row, 1).Value = $sender
        $ws.Cells.Item($row, 2).Value = "Angpt1"
        $ws.Cells.Item($row, 3).Value = "Itgb1"
        $ws.Cells.Item($row, 4).Value = $target

        $senderVals = $bySender[$sender]
        $ws.Cells.Item($row, 5).Value  = $senderVals[0]
        $ws.Cells.Item($row, 6).Value  = $senderVals[1]
        $ws.Cells.Item($row, 7).Value  = $senderVals[2]
        $ws.Cells.Item($row, 8).Value  = $senderVals[3]
        $ws.Cells.Item($row, 9).Value  = $senderVals[4]
        $ws.Cells.Item($row, 10).Value = $senderVals[5]

        $targetVals = $byTarget[$target]
        $ws.Cells.Item($row, 11).Value = $targetVals[0]
        $ws.Cells.Item($row, 12).Value = $targetVals[1]
        $ws.Cells.Item($row, 13).Value = $targetVals[2]
        $ws.Cells.Item($row, 14).Value = $targetVals[3]
        $ws.Cells.Item($row, 15).Value = $targetVals[4]
        $ws.Cells.Item($row, 16).Value = $targetVals[5]

        $pairVals = $byPair["$sender|$target"]
        $ws.Cells.Item($row, 17).Value = $pairVals[0]
        $ws.Cells.Item($row, 18).Value = $pairVals[1]
        $ws.Cells.Item($row, 19).Value = $pairVals[2]
        $ws.Cells.Item($row, 20).Value = $pairVals[3]

        $row = $row + 1
    }
}
